# RSTK-9519 - Disassembly Work Order: optimize WO test cases.
#
# Summary of the change:
#  - "Create Disassembly WO" sheet: the second disassembly-item test row is
#    updated to use the "Serial track" item / qty 2, and the now-redundant
#    third test row is removed entirely.
#  - "Consumable Components" sheet: the B8 component reference is switched to
#    "Pro-child4 (NO Track)" and the two trailing test rows (9 and 10) that
#    duplicated this scenario are removed.
#  - The active sheet / selection moves from "Create Disassembly WO" to
#    "Consumable Components", and the selection on "Derived Components" is
#    updated as well (view-state only, no data change there).

$wb = $excel.ActiveWorkbook

# --- Sheet "Create Disassembly WO" -----------------------------------------
$ws1 = $wb.Worksheets.Item("Create Disassembly WO")

# Row 2 now describes the "Pro-Disassembley Serial (Serial track)" item with
# qty 2 (previously qty 3 against the lot+serial tracked item), and no longer
# carries the "Add consumable components?" flag in column F.
$ws1.Range("A2").Value = "Pro-Disassembley Serial (Serial track)"
$ws1.Range("B2").Value = 2
$ws1.Range("F2").ClearContents() | Out-Null

# Former row 3 (the duplicate "Serial track" qty-2 test case) is removed, its
# data having been folded into row 2 above.
$ws1.Rows(3).Delete() | Out-Null

# --- Sheet "Derived Components" ---------------------------------------------
$ws2 = $wb.Worksheets.Item("Derived Components")

# View-state only: selection moves onto rows 7:8.
$ws2.Rows("7:8").Select() | Out-Null

# --- Sheet "Consumable Components" ------------------------------------------
$ws3 = $wb.Worksheets.Item("Consumable Components")

# Row 8's derived component reference switches to "Pro-child4 (NO Track)".
$ws3.Range("B8").Value = "Pro-child4 (NO Track)"

# Rows 9 and 10 (the other two consumable-component scenarios) are removed.
$ws3.Rows("9:10").Delete() | Out-Null

# This sheet becomes the active tab/selected sheet, with the selection
# resting on the now-empty rows 9:10.
$ws3.Rows("9:10").Select() | Out-Null
$ws3.Activate() | Out-Null
